$d = $word.ActiveDocument

# The document lists many "childless tag" paragraphs (e.g. "PUMP:HRD:0000",
# "PUMP:HTR:100", ..., "PUMPHTR:200") right after the intro paragraph
# ("These are the childless tags that were found in the documents: ").
# The edit removes all of those listing paragraphs, leaving only the
# title paragraph and the intro paragraph before the section properties.

$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("PUMP")) {
        $startIndex = $i
        break
    }
}

if ($startIndex -ge 1) {
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
